# Regenerate the localization handoff report:
#   - The source file was renamed/regenerated under a new GUID
#     (855bebd0-e14f-4b76-afd8-bfd13c3e8764 -> 92fdfdda-9fd8-4f38-bbf3-e1cc1a5b84db)
#   - A fresh handoff round was produced (new xlf package hash, new timestamps)
#   - Neither locale has been handed back yet, so the "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" columns reset to
#     "not yet" values on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "855bebd0-e14f-4b76-afd8-bfd13c3e8764"
$newGuid = "92fdfdda-9fd8-4f38-bbf3-e1cc1a5b84db"
$newHash = "de93b67dcf5156182d731a970b015bb6578e7e2c"

$newFileName      = "$newGuid.md"
$newPathAndName   = "e2e\$newGuid.md"
$newGenerateDate  = "2016-08-26 04:57:21"
$newZhXlf         = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf         = "$newGuid.$newHash.de-de.xlf"
$newZhHandoffDate = "2016-08-26 04:57:17"
$neverHandedBack  = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newFileName
$ov.Range("G2").Value = $newGenerateDate

# Re-point the B2 hyperlink's display text at the new path (keep same target
# URL - it is unchanged in the diff).
$ovUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8a8df3dba688edf6ac190d4667facb86d47ce21/e2e/$oldGuid.md"
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ovUrl, "", "", $newPathAndName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newFileName
$zh.Range("G2").Value = $newZhXlf
$zh.Range("H2").Value = $newZhHandoffDate
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("J2").Style = "Normal"
$zh.Range("K2").Value = $neverHandedBack

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8a8df3dba688edf6ac190d4667facb86d47ce21/e2e/$oldGuid.md"
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhUrl, "", "", $newFileName)

$zh.Columns.Item(9).AutoFit()
$zh.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newFileName
$de.Range("G2").Value = $newDeXlf
$de.Range("H2").Value = $newGenerateDate
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("J2").Style = "Normal"
$de.Range("K2").Value = $neverHandedBack

$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8a8df3dba688edf6ac190d4667facb86d47ce21/e2e/$oldGuid.md"
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deUrl, "", "", $newFileName)

$de.Columns.Item(9).AutoFit()
$de.Columns.Item(10).AutoFit()
